# Add "Varun" to the list of people assigned to use cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "Varun"
$ws.Range("B4").Value = "Varun"
$ws.Range("B10").Value = "Varun"

$ws.Range("B4").Select()
